$d = $word.ActiveDocument

# Each replacement targets a unique p-value cell in the Mantel correlogram table.
# Find.Execute params: FindText, MatchCase, MatchWholeWord, MatchWildcards,
# MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace

$d.Content.Find.Execute("0.433", $true, $true, $false, $false, $false,
                         $true, 1, $false, "0.431", 2)

$d.Content.Find.Execute("0.478", $true, $true, $false, $false, $false,
                         $true, 1, $false, "0.444", 2)

$d.Content.Find.Execute("0.575", $true, $true, $false, $false, $false,
                         $true, 1, $false, "0.59", 2)

$d.Content.Find.Execute("0.767", $true, $true, $false, $false, $false,
                         $true, 1, $false, "0.787", 2)

$d.Content.Find.Execute("0.959", $true, $true, $false, $false, $false,
                         $true, 1, $false, "0.984", 2)
